# Applies the commit "Atualizacao de bases das ligas, do dia: 28-05-2024 as 20:56"
# Swaps several match rows (same-date fixtures reordered) and updates odds
# for a handful of upcoming fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps: rows sharing the same match date were reordered. ---
# For each cycle of rows, snapshot the B:AD payload of every row in the
# cycle first, then write it back rotated so no data is lost mid-way.

# Cycle: 22, 23
$row22 = $ws.Range("B22:AD22").Value2
$row23 = $ws.Range("B23:AD23").Value2
$ws.Range("B22:AD22").Value2 = $row23
$ws.Range("B23:AD23").Value2 = $row22

# Cycle: 24, 25
$row24 = $ws.Range("B24:AD24").Value2
$row25 = $ws.Range("B25:AD25").Value2
$ws.Range("B24:AD24").Value2 = $row25
$ws.Range("B25:AD25").Value2 = $row24

# Cycle: 29, 30
$row29 = $ws.Range("B29:AD29").Value2
$row30 = $ws.Range("B30:AD30").Value2
$ws.Range("B29:AD29").Value2 = $row30
$ws.Range("B30:AD30").Value2 = $row29

# Cycle: 35, 38
$row35 = $ws.Range("B35:AD35").Value2
$row38 = $ws.Range("B38:AD38").Value2
$ws.Range("B35:AD35").Value2 = $row38
$ws.Range("B38:AD38").Value2 = $row35

# Cycle: 36, 37
$row36 = $ws.Range("B36:AD36").Value2
$row37 = $ws.Range("B37:AD37").Value2
$ws.Range("B36:AD36").Value2 = $row37
$ws.Range("B37:AD37").Value2 = $row36

# Cycle: 39, 41, 40
$row39 = $ws.Range("B39:AD39").Value2
$row41 = $ws.Range("B41:AD41").Value2
$row40 = $ws.Range("B40:AD40").Value2
$ws.Range("B39:AD39").Value2 = $row41
$ws.Range("B41:AD41").Value2 = $row40
$ws.Range("B40:AD40").Value2 = $row39

# Cycle: 43, 45, 44
$row43 = $ws.Range("B43:AD43").Value2
$row45 = $ws.Range("B45:AD45").Value2
$row44 = $ws.Range("B44:AD44").Value2
$ws.Range("B43:AD43").Value2 = $row45
$ws.Range("B45:AD45").Value2 = $row44
$ws.Range("B44:AD44").Value2 = $row43

# Cycle: 56, 57, 58, 59
$row56 = $ws.Range("B56:AD56").Value2
$row57 = $ws.Range("B57:AD57").Value2
$row58 = $ws.Range("B58:AD58").Value2
$row59 = $ws.Range("B59:AD59").Value2
$ws.Range("B56:AD56").Value2 = $row57
$ws.Range("B57:AD57").Value2 = $row58
$ws.Range("B58:AD58").Value2 = $row59
$ws.Range("B59:AD59").Value2 = $row56

# Cycle: 154, 155, 157, 156
$row154 = $ws.Range("B154:AD154").Value2
$row155 = $ws.Range("B155:AD155").Value2
$row157 = $ws.Range("B157:AD157").Value2
$row156 = $ws.Range("B156:AD156").Value2
$ws.Range("B154:AD154").Value2 = $row155
$ws.Range("B155:AD155").Value2 = $row157
$ws.Range("B157:AD157").Value2 = $row156
$ws.Range("B156:AD156").Value2 = $row154

# Cycle: 160, 161
$row160 = $ws.Range("B160:AD160").Value2
$row161 = $ws.Range("B161:AD161").Value2
$ws.Range("B160:AD160").Value2 = $row161
$ws.Range("B161:AD161").Value2 = $row160

# Cycle: 164, 169
$row164 = $ws.Range("B164:AD164").Value2
$row169 = $ws.Range("B169:AD169").Value2
$ws.Range("B164:AD164").Value2 = $row169
$ws.Range("B169:AD169").Value2 = $row164

# Cycle: 193, 196
$row193 = $ws.Range("B193:AD193").Value2
$row196 = $ws.Range("B196:AD196").Value2
$ws.Range("B193:AD193").Value2 = $row196
$ws.Range("B196:AD196").Value2 = $row193

# Cycle: 195, 197
$row195 = $ws.Range("B195:AD195").Value2
$row197 = $ws.Range("B197:AD197").Value2
$ws.Range("B195:AD195").Value2 = $row197
$ws.Range("B197:AD197").Value2 = $row195

# --- Odds-only updates for upcoming (not-yet-played) fixtures. ---
$ws.Range("O239").Value2 = 1.666
$ws.Range("P239").Value2 = 3.8
$ws.Range("Q239").Value2 = 4.75
$ws.Range("R239").Value2 = -0.75
$ws.Range("S239").Value2 = 1.875
$ws.Range("T239").Value2 = 1.975
$ws.Range("U239").Value2 = 2.75
$ws.Range("V239").Value2 = 1.825
$ws.Range("W239").Value2 = 2.025

$ws.Range("O240").Value2 = 1.8
$ws.Range("P240").Value2 = 3.9
$ws.Range("Q240").Value2 = 3.8
$ws.Range("R240").Value2 = -0.5
$ws.Range("S240").Value2 = 1.825
$ws.Range("T240").Value2 = 2.025
$ws.Range("U240").Value2 = 2.75
$ws.Range("V240").Value2 = 1.825
$ws.Range("W240").Value2 = 2.025

$ws.Range("O241").Value2 = 1.7
$ws.Range("P241").Value2 = 4
$ws.Range("Q241").Value2 = 4.2
$ws.Range("R241").Value2 = -0.75
$ws.Range("S241").Value2 = 1.925
$ws.Range("T241").Value2 = 1.925
$ws.Range("U241").Value2 = 3
$ws.Range("V241").Value2 = 1.875
$ws.Range("W241").Value2 = 1.975

$ws.Range("O242").Value2 = 3.2
$ws.Range("P242").Value2 = 3.6
$ws.Range("Q242").Value2 = 2.1
$ws.Range("R242").Value2 = 0.25
$ws.Range("S242").Value2 = 2
$ws.Range("T242").Value2 = 1.85
$ws.Range("U242").Value2 = 3
$ws.Range("V242").Value2 = 1.975
$ws.Range("W242").Value2 = 1.875

$ws.Range("O243").Value2 = 1.727
$ws.Range("P243").Value2 = 3.9
$ws.Range("Q243").Value2 = 4.2
$ws.Range("R243").Value2 = -0.75
$ws.Range("S243").Value2 = 1.975
$ws.Range("T243").Value2 = 1.875
$ws.Range("U243").Value2 = 2.75
$ws.Range("V243").Value2 = 1.825
$ws.Range("W243").Value2 = 2.025

$ws.Range("O244").Value2 = 2.1
$ws.Range("P244").Value2 = 3.75
$ws.Range("Q244").Value2 = 3.1
$ws.Range("R244").Value2 = -0.25
$ws.Range("S244").Value2 = 1.85
$ws.Range("T244").Value2 = 2
$ws.Range("U244").Value2 = 2.75
$ws.Range("V244").Value2 = 1.85
$ws.Range("W244").Value2 = 2

$ws.Range("O245").Value2 = 2.05
$ws.Range("P245").Value2 = 3.8
$ws.Range("Q245").Value2 = 3
$ws.Range("R245").Value2 = -0.25
$ws.Range("S245").Value2 = 1.85
$ws.Range("T245").Value2 = 2
$ws.Range("U245").Value2 = 3
$ws.Range("V245").Value2 = 1.9
$ws.Range("W245").Value2 = 1.95
